$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Damian Lillard"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Milwaukee Bucks"
$ws.Range("A3").Value = "Derrick White"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Boston Celtics"
$ws.Range("A4").Value = "Jared McCain"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Philadelphia 76ers"
$ws.Range("A5").Value = "Cam Thomas"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Brooklyn Nets"
$ws.Range("A6").Value = "Brandon Ingram"
$ws.Range("B6").Value = "SG,SF,PF"
$ws.Range("C6").Value = "New Orleans Pelicans"
$ws.Range("A7").Value = "Cameron Johnson"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Brooklyn Nets"
$ws.Range("A8").Value = "Yves Missi"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "New Orleans Pelicans"
$ws.Range("A9").Value = "Robert Williams III"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Portland Trail Blazers"
$ws.Range("A10").Value = "Bam Adebayo"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Miami Heat"
$ws.Range("A11").Value = "Dalton Knecht"
$ws.Range("B11").Value = "SG"
$ws.Range("C11").Value = "Los Angeles Lakers"
$ws.Range("A12").Value = "Anthony Davis"
$ws.Range("B12").Value = "PF,C"
$ws.Range("C12").Value = "Los Angeles Lakers"
$ws.Range("A13").Value = "Julius Randle"
$ws.Range("B13").Value = "PF"
$ws.Range("C13").Value = "Minnesota Timberwolves"
$ws.Range("A14").Value = "LaMelo Ball"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Charlotte Hornets"
$ws.Range("A15").Value = "Brandon Miller"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Charlotte Hornets"
$ws.Range("A16").Value = "Isaiah Hartenstein"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Oklahoma City Thunder"
$ws.Range("A17").Value = "Zion Williamson"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "New Orleans Pelicans"
$ws.Range("A18").Value = "Cade Cunningham"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Detroit Pistons"

# Remove the now-unused last row (old data had 18 rows of players, new data has 17)
$ws.Rows("19").Delete()
